# Applies corrected Diebold-Mariano DM_Stat (col C) and P_Value (col D)
# figures for rows 2-11 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = 0.3479797917628265;  D = 0.7311629094229573 },
    @{ Row = 3;  C = 0.9463228493038347;  D = 0.3542628728215815 },
    @{ Row = 4;  C = 0.6213703784046531;  D = 0.540741644997051  },
    @{ Row = 5;  C = 0.2075958078937248;  D = 0.8374542320256557 },
    @{ Row = 6;  C = 0.5004385915664157;  D = 0.6217357592130806 },
    @{ Row = 7;  C = 0.3122318775782317;  D = 0.7578044416166161 },
    @{ Row = 8;  C = -0.1138489236722984; D = 0.9103901511857988 },
    @{ Row = 9;  C = -0.1154374038163631; D = 0.9091456207495248 },
    @{ Row = 10; C = -0.7804114403766289; D = 0.4434650772195803 },
    @{ Row = 11; C = -0.4359494125534956; D = 0.6671212996994598 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
